$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.895.95"
$ws.Range("E2").Value = "  -0.42%  "

$ws.Range("D3").Value = "1.632.21"
$ws.Range("E3").Value = "  -0.94%  "

$ws.Range("D5").Value = "'211.73"
$ws.Range("E5").Value = "  -0.92%  "

$ws.Range("D6").Value = "'0.523"
$ws.Range("E6").Value = "  -0.90%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "'23.22"
$ws.Range("E8").Value = "  -1.06%  "

$ws.Range("E9").Value = "  -3.20%  "

$ws.Range("E10").Value = "  -0.38%  "

$ws.Range("E11").Value = "  +0.73%  "

$ws.Range("D12").Value = "1.863.87"
$ws.Range("E12").Value = "  -0.94%  "

$ws.Range("D13").Value = "1.642.96"
$ws.Range("E13").Value = "  -0.28%  "

$ws.Range("E14").Value = "  -0.50%  "

$ws.Range("E15").Value = "  -0.03%  "

$ws.Range("D16").Value = "'65.22"
$ws.Range("E16").Value = "  -0.60%  "

$ws.Range("D17").Value = "27.898.47"
$ws.Range("E17").Value = "  -0.38%  "

$ws.Range("D18").Value = "'229.91"
$ws.Range("E18").Value = "  -1.39%  "

$ws.Range("E20").Value = "  -2.42%  "

$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("E23").Value = "  -2.96%  "

$ws.Range("E24").Value = "  -4.26%  "

$ws.Range("D25").Value = "'153.82"
$ws.Range("E25").Value = "  +0.79%  "

$ws.Range("E26").Value = "  +0.35%  "

$ws.Range("E27").Value = "  -0.76%  "

$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("D30").Value = "'1.18"
$ws.Range("E30").Value = "  -0.97%  "

$ws.Range("D31").Value = "'0.0482"
$ws.Range("E31").Value = "  -0.28%  "

$ws.Range("D32").Value = "'3.37"
$ws.Range("E32").Value = "  +0.41%  "

$ws.Range("E33").Value = "  -0.91%  "

$ws.Range("D34").Value = "1.397.82"
$ws.Range("E34").Value = "  -3.55%  "

$ws.Range("E35").Value = "  -0.39%  "

$ws.Range("E36").Value = "  +9.73%  "

$ws.Range("E37").Value = "  +1.34%  "

$ws.Range("E38").Value = "  +0.44%  "

$ws.Range("E39").Value = "  -0.36%  "

$ws.Range("E40").Value = "  -2.14%  "

$ws.Range("E41").Value = "  -0.21%  "

$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("D43").Value = "'66.85"
$ws.Range("E43").Value = "  -3.80%  "

$ws.Range("E44").Value = "  +2.57%  "

$ws.Range("E45").Value = "  +0.95%  "

$ws.Range("E46").Value = "  -1.09%  "

$ws.Range("D47").Value = "1.773.49"
$ws.Range("E47").Value = "  -0.98%  "

$ws.Range("D48").Value = "'87.61"
$ws.Range("E48").Value = "  -1.61%  "

$ws.Range("D49").Value = "0.0₆0105"
$ws.Range("E49").Value = "  +0.31%  "

$ws.Range("E50").Value = "  -0.91%  "
